# EquivalencesPatrons.xlsx update:
# - The lab's "Command" pattern row now maps to "PerspectiveCommand"
#   instead of "ImageCommand" (last diagram + skeleton changed).
# - The whole equivalence table (D6:E19) gets a thin box border around
#   every cell.
# - The sheet's saved selection now highlights the whole table
#   (D6:E19) with D6 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell text that referenced the old lab element name.
$ws.Range("E14").Value = "PerspectiveCommand"

# Add a thin border around every cell of the equivalence table.
$tableRange = $ws.Range("D6:E19")
$tableRange.Borders.LineStyle = 1

# Select the whole table, with D6 as the active cell, matching the
# selection that gets persisted in the sheet view.
$tableRange.Select() | Out-Null
